$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "period" column (G) values: replace the verbose labels
# ("Evening (6pm - 12am)" / "Night (12am - 6am)") with short lowercase tags
$ws.Range("G2").Value = "evening"
$ws.Range("G3").Value = "evening"
$ws.Range("G4").Value = "evening"
$ws.Range("G5").Value = "evening"
$ws.Range("G6").Value = "evening"
$ws.Range("G7").Value = "night"

# Add new header columns "song" and "artist"
$ws.Range("J1").Value = "song"
$ws.Range("K1").Value = "artist"

# Split the "song-artist" (H) column into "song" (J) and "artist" (K) on " - "
$songArtist = @(
    "Lacrimosa - Mozart",
    "Serenade - Schubert",
    "Mad About You - Hooverphonic",
    "Wasting My Young Years - London Grammar",
    "Dusk Till Dawn - ZAYN ft. Sia",
    "Honey Bee - Madrugada "
)

for ($i = 0; $i -lt $songArtist.Length; $i++) {
    $row = $i + 2
    $parts = $songArtist[$i] -split ' - ', 2
    $ws.Range("J$row").Value = $parts[0] + " "
    $ws.Range("K$row").Value = " " + $parts[1]
}

# Set the column widths for the new columns, matching bestFit sizing
$ws.Columns.Item(10).ColumnWidth = 21.6640625
$ws.Columns.Item(11).ColumnWidth = 15.77734375

# Update the active selection to K3, matching the final state
$ws.Range("K3").Select()
